$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Label" column header (H1) and per-row binary label (H2:H21)
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Row 2: Control 0
$ws.Range("D2").Value = 0.4373152328734827
$ws.Range("E2").Value = 0.4373152328734827
$ws.Range("H2").Value = 0

# Row 3: Control 6
$ws.Range("D3").Value = 0.5695225106817564
$ws.Range("E3").Value = 0.5695225106817564
$ws.Range("H3").Value = 0

# Row 4: Control 9
$ws.Range("D4").Value = 0.6790430128590882
$ws.Range("E4").Value = 0.6790430128590882
$ws.Range("H4").Value = 0

# Row 5: Control 24
$ws.Range("D5").Value = 0.5318985243300074
$ws.Range("E5").Value = 0.5318985243300074
$ws.Range("H5").Value = 0

# Row 6: Control 32
$ws.Range("D6").Value = 0.5060624540390143
$ws.Range("E6").Value = 0.5060624540390143
$ws.Range("H6").Value = 0

# Row 7: MDD 27
$ws.Range("D7").Value = 0.5617026332579738
$ws.Range("E7").Value = 0.4382973667420262
$ws.Range("H7").Value = 1

# Row 8: MDD 47
$ws.Range("D8").Value = 0.4849879024211191
$ws.Range("E8").Value = 0.515012097578881
$ws.Range("H8").Value = 1

# Row 9: MDD 13
$ws.Range("D9").Value = 0.4872274281999144
$ws.Range("E9").Value = 0.5127725718000855
$ws.Range("H9").Value = 1

# Row 10: MDD 25
$ws.Range("D10").Value = 0.5074442587366423
$ws.Range("E10").Value = 0.4925557412633577
$ws.Range("H10").Value = 1

# Row 11: MDD 5
$ws.Range("D11").Value = 0.7989714679607819
$ws.Range("E11").Value = 0.2010285320392181
$ws.Range("F11").Value = 0.6940997242927551
$ws.Range("H11").Value = 1

# Row 12: Control 0
$ws.Range("H12").Value = 0

# Row 13: Control 6
$ws.Range("H13").Value = 0

# Row 14: Control 9
$ws.Range("H14").Value = 0

# Row 15: Control 24
$ws.Range("H15").Value = 0

# Row 16: Control 32
$ws.Range("H16").Value = 0

# Row 17: MDD 27
$ws.Range("H17").Value = 1

# Row 18: MDD 47
$ws.Range("H18").Value = 1

# Row 19: MDD 13
$ws.Range("H19").Value = 1

# Row 20: MDD 25
$ws.Range("H20").Value = 1

# Row 21: MDD 5
$ws.Range("H21").Value = 1

